$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting from the last existing data row (14) down to the new rows (15-23) ---
$ws.Range("A14:P14").Copy() | Out-Null
$ws.Range("A15:P23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 1: set string cells in the precise order required so that new shared strings ---
# --- are appended to sharedStrings.xml in index order 23..45 matching the target diff. ---
$ws.Range("A15").Value = "Stream-crude glycerol"
$ws.Range("A16").Value = "Stream-pure glycerine"
$ws.Range("A17").Value = "Stream-cellulase"
$ws.Range("A19").Value = "Pretreatment reactor system"
$ws.Range("A20").Value = "Pretreatment and saccharification"
$ws.Range("A22").Value = "Cofermenation"
$ws.Range("B4").Value = "Cane lipid content [dry wt. %]"
$ws.Range("B5").Value = "Relative sorghum lipid content [dry wt. %]"
$ws.Range("B6").Value = "Lipid retention [%]"
$ws.Range("B7").Value = "Bagasse lipid extraction efficiency [%]"
$ws.Range("B8").Value = "Capacity [ton/hr]"
$ws.Range("B9").Value = "Price [USD/gal]"
$ws.Range("B11").Value = "Price [USD/cf]"
$ws.Range("B12").Value = "Electricity price [USD/kWh]"
$ws.Range("B13").Value = "Operating days [day/yr]"
$ws.Range("B14").Value = "IRR [%]"
$ws.Range("B15").Value = "Price [USD/kg]"
$ws.Range("B18").Value = "Cellulase loading [wt. % cellulose]"
$ws.Range("B19").Value = "Base cost [million USD]"
$ws.Range("B20").Value = "Glucose yield [%]"
$ws.Range("B21").Value = "Xylose yield [%]"
$ws.Range("B22").Value = "Glucose to ethanol yield [%]"
$ws.Range("B23").Value = "Xylose to ethanol yield [%]"

# --- Step 2: set remaining string cells (values already exist in sharedStrings) ---
$ws.Range("A4").Value = "Stream-lipidcane"
$ws.Range("A9").Value = "Stream-ethanol"
$ws.Range("A10").Value = "Stream-biodiesel"
$ws.Range("B10").Value = "Price [USD/gal]"
$ws.Range("A11").Value = "Stream-natural gas"
$ws.Range("A12").Value = "biorefinery"
$ws.Range("B16").Value = "Price [USD/kg]"
$ws.Range("B17").Value = "Price [USD/kg]"

# --- Step 3: set numeric data cells ---
$ws.Range("C4").Value = 0.2469923450796938
$ws.Range("D4").Value = 0.986758600910344
$ws.Range("E4").Value = -0.9064206340488252
$ws.Range("G4").Value = -0.2607400958056038
$ws.Range("H4").Value = 0.14821650026466
$ws.Range("I4").Value = -0.01894716392588655
$ws.Range("J4").Value = 0.09486322049384942
$ws.Range("K4").Value = -0.05555093415803736
$ws.Range("L4").Value = -0.01780926090437044
$ws.Range("M4").Value = -0.03925711030628441
$ws.Range("O4").Value = 0.1799079849563194
$ws.Range("P4").Value = -0.0597553403422136
$ws.Range("C5").Value = 0.002719967628798705
$ws.Range("D5").Value = -0.01263442082537683
$ws.Range("E5").Value = 0.003635610673424427
$ws.Range("G5").Value = -0.01400248577609943
$ws.Range("H5").Value = 0.002204892856195714
$ws.Range("I5").Value = -0.00162593900903756
$ws.Range("J5").Value = 0.006856403604230787
$ws.Range("K5").Value = 0.02680977630439105
$ws.Range("L5").Value = -0.01568699045147962
$ws.Range("M5").Value = 0.001161293422451737
$ws.Range("O5").Value = 0.001493080475723219
$ws.Range("P5").Value = 0.0041438075257523
$ws.Range("C6").Value = -0.005504406844176273
$ws.Range("D6").Value = 0.001741205157648206
$ws.Range("E6").Value = -0.01351631305265252
$ws.Range("G6").Value = -0.02795857897434316
$ws.Range("H6").Value = 0.008176599399063974
$ws.Range("I6").Value = -0.002647487913899516
$ws.Range("J6").Value = 0.00002309856093059482
$ws.Range("K6").Value = -0.03712544558101782
$ws.Range("L6").Value = -0.1142591916103677
$ws.Range("M6").Value = -0.005349988629999545
$ws.Range("O6").Value = -0.02427687293907492
$ws.Range("P6").Value = 0.01058900922356037
$ws.Range("C7").Value = 0.07469037821961512
$ws.Range("D7").Value = 0.1372498183859927
$ws.Range("E7").Value = 0.01711218903648756
$ws.Range("G7").Value = 0.2961772888710915
$ws.Range("H7").Value = 0.01034132422165297
$ws.Range("I7").Value = -0.0005958812398352496
$ws.Range("J7").Value = -0.0002316977853346404
$ws.Range("K7").Value = 0.1568886093475444
$ws.Range("L7").Value = 0.9900151562406062
$ws.Range("M7").Value = -0.02128182440327298
$ws.Range("O7").Value = 0.2410137919125517
$ws.Range("P7").Value = 0.0379576852623074
$ws.Range("C8").Value = 0.2129644699905788
$ws.Range("D8").Value = -0.01777513271100531
$ws.Range("E8").Value = 0.02097542474301698
$ws.Range("G8").Value = -0.7436686520507461
$ws.Range("H8").Value = 0.9552592068503682
$ws.Range("I8").Value = 0.9576457947858317
$ws.Range("J8").Value = 0.0401082729918821
$ws.Range("K8").Value = 0.03122272915290916
$ws.Range("L8").Value = 0.006746925197877007
$ws.Range("M8").Value = 0.0924876377795055
$ws.Range("O8").Value = 0.01494019326960773
$ws.Range("P8").Value = 0.3387088714683548
$ws.Range("C9").Value = 0.736641527993661
$ws.Range("D9").Value = -0.01412434578097383
$ws.Range("E9").Value = 0.01779446211977848
$ws.Range("G9").Value = 0.02516545819061832
$ws.Range("H9").Value = -0.02771862523674501
$ws.Range("I9").Value = -0.02565199811407992
$ws.Range("J9").Value = -0.005749891191651616
$ws.Range("K9").Value = -0.4006174151926966
$ws.Range("L9").Value = 0.001468714522748581
$ws.Range("M9").Value = 0.01017482709499308
$ws.Range("O9").Value = 0.002857585362303414
$ws.Range("P9").Value = 0.0187662575826503
$ws.Range("C10").Value = 0.4041355319414213
$ws.Range("D10").Value = -0.01573430683737227
$ws.Range("E10").Value = 0.00729812150792486
$ws.Range("G10").Value = -0.01160250488010019
$ws.Range("H10").Value = 0.00920472305618892
$ws.Range("I10").Value = 0.007491558251662329
$ws.Range("J10").Value = 0.01244500749338446
$ws.Range("K10").Value = 0.863551453326058
$ws.Range("L10").Value = -0.01375502282220091
$ws.Range("M10").Value = -0.002713292556531702
$ws.Range("O10").Value = 0.0250781325551253
$ws.Range("P10").Value = -0.01135464381418575
$ws.Range("C11").Value = -0.03624571911382876
$ws.Range("D11").Value = 0.01644022884960915
$ws.Range("E11").Value = -0.01249146213165848
$ws.Range("G11").Value = -0.002402886624115465
$ws.Range("H11").Value = 0.01047237325089493
$ws.Range("I11").Value = 0.008194044615761783
$ws.Range("J11").Value = 0.001916086253195283
$ws.Range("K11").Value = 0.01308897767555911
$ws.Range("L11").Value = -0.0006360300734412029
$ws.Range("M11").Value = 0.001527963325118533
$ws.Range("O11").Value = 0.01055817594232704
$ws.Range("P11").Value = -0.0128029449921178
$ws.Range("C12").Value = 0.01126975485079019
$ws.Range("D12").Value = -0.01090531579621263
$ws.Range("E12").Value = 0.01224867610594704
$ws.Range("G12").Value = -0.00818205104728204
$ws.Range("H12").Value = 0.01806721982668879
$ws.Range("I12").Value = 0.02000939523237581
$ws.Range("J12").Value = 0.007664973716806461
$ws.Range("K12").Value = -0.02390648783625951
$ws.Range("L12").Value = 0.003811060760442429
$ws.Range("M12").Value = 0.006088397427535897
$ws.Range("O12").Value = 0.00836607979064319
$ws.Range("P12").Value = 0.003910724508428979
$ws.Range("C13").Value = 0.1085119401804776
$ws.Range("D13").Value = -0.001362051126482045
$ws.Range("E13").Value = -0.0006576818183072727
$ws.Range("G13").Value = 0.001771656166866246
$ws.Range("H13").Value = -0.009877623851104952
$ws.Range("I13").Value = 0.267955776926231
$ws.Range("J13").Value = -0.0016910668041297
$ws.Range("K13").Value = 0.003968628926745156
$ws.Range("L13").Value = -0.01549295409171816
$ws.Range("M13").Value = 0.001467111802684472
$ws.Range("O13").Value = -0.0195405726856229
$ws.Range("P13").Value = -0.003022912632916505
$ws.Range("C14").Value = -0.2306115292884611
$ws.Range("D14").Value = 0.01709096890763875
$ws.Range("E14").Value = -0.01556795620671825
$ws.Range("G14").Value = -0.01932036298081452
$ws.Range("H14").Value = 0.02125225736209029
$ws.Range("I14").Value = 0.02607368293094731
$ws.Range("J14").Value = 0.01794493694696562
$ws.Range("K14").Value = -0.02285042462601699
$ws.Range("L14").Value = 0.006835094577403783
$ws.Range("M14").Value = -0.003613125552525021
$ws.Range("O14").Value = 0.01545468714618748
$ws.Range("P14").Value = 0.009024059592962382
$ws.Range("C15").Value = 0.006306360444254418
$ws.Range("D15").Value = -0.009873999466959977
$ws.Range("E15").Value = 0.01306457313058292
$ws.Range("G15").Value = 0.007759988566399542
$ws.Range("H15").Value = 0.0003115197244607889
$ws.Range("I15").Value = 0.004215266088610644
$ws.Range("J15").Value = -0.03462567221099908
$ws.Range("K15").Value = 0.03986642444265697
$ws.Range("L15").Value = 0.003834531129381244
$ws.Range("M15").Value = -0.00532308203692328
$ws.Range("O15").Value = 0.001035289097411564
$ws.Range("P15").Value = -0.01999254799970192
$ws.Range("C16").Value = -0.03460892144835685
$ws.Range("D16").Value = 0.008261616714464669
$ws.Range("E16").Value = -0.007483931723357268
$ws.Range("G16").Value = -0.02657627127105085
$ws.Range("H16").Value = 0.0306785727151429
$ws.Range("I16").Value = 0.02145504949820198
$ws.Range("J16").Value = -0.0004381455536520081
$ws.Range("K16").Value = -0.05574883835795352
$ws.Range("L16").Value = -0.0005849899433995977
$ws.Range("M16").Value = -0.00319263199970528
$ws.Range("O16").Value = 0.005045688393827534
$ws.Range("P16").Value = 0.01187203602688144
$ws.Range("C17").Value = -0.04764594046583761
$ws.Range("D17").Value = 0.03858327696733107
$ws.Range("E17").Value = -0.03801147377645895
$ws.Range("G17").Value = -0.0314825474993019
$ws.Range("H17").Value = 0.05101708322468332
$ws.Range("I17").Value = 0.04796831395073255
$ws.Range("J17").Value = 0.0264738707145793
$ws.Range("K17").Value = -0.03717406689496267
$ws.Range("L17").Value = -0.003088877787555111
$ws.Range("M17").Value = -0.004457309266292371
$ws.Range("O17").Value = -0.009309064404362574
$ws.Range("P17").Value = 0.02608625422745017
$ws.Range("C18").Value = 0.03419123464764938
$ws.Range("D18").Value = 0.003723822580952902
$ws.Range("E18").Value = -0.007510390380415614
$ws.Range("G18").Value = -0.008686688795467551
$ws.Range("H18").Value = 0.01364614153784566
$ws.Range("I18").Value = 0.01683920947356838
$ws.Range("J18").Value = -0.01131136955971246
$ws.Range("K18").Value = 0.007176040319041612
$ws.Range("L18").Value = 0.01595050767802031
$ws.Range("M18").Value = -0.005984205359368213
$ws.Range("O18").Value = -0.005041146345645853
$ws.Range("P18").Value = 0.006408139168325566
$ws.Range("C19").Value = -0.05760135849605433
$ws.Range("D19").Value = -0.0155791351191654
$ws.Range("E19").Value = 0.02634056562962262
$ws.Range("G19").Value = 0.02074488927779557
$ws.Range("H19").Value = 0.2312987370599495
$ws.Range("I19").Value = 0.004456985650279426
$ws.Range("J19").Value = -0.01245136379764054
$ws.Range("K19").Value = -0.01846022550640902
$ws.Range("L19").Value = 0.01639494535979781
$ws.Range("M19").Value = -0.01441665849666634
$ws.Range("O19").Value = -0.001334751605390064
$ws.Range("P19").Value = 0.1897123344044933
$ws.Range("C20").Value = 0.03539219229568769
$ws.Range("D20").Value = -0.003172570398902815
$ws.Range("E20").Value = 0.1627699066387962
$ws.Range("G20").Value = 0.1853130940045237
$ws.Range("H20").Value = 0.009657139874285594
$ws.Range("I20").Value = -0.001009648360385934
$ws.Range("J20").Value = 0.01215310877762445
$ws.Range("K20").Value = 0.0004549330741973229
$ws.Range("L20").Value = -0.01193075970923039
$ws.Range("M20").Value = 0.3096087550083502
$ws.Range("O20").Value = 0.05643515707340628
$ws.Range("P20").Value = 0.007276784739071388
$ws.Range("C21").Value = 0.1094825405393016
$ws.Range("D21").Value = 0.01979661909586476
$ws.Range("E21").Value = 0.2308598804343952
$ws.Range("G21").Value = 0.06036057563042301
$ws.Range("H21").Value = -0.001099089451963578
$ws.Range("I21").Value = 0.02411104560444182
$ws.Range("J21").Value = -0.04211265822463477
$ws.Range("K21").Value = 0.03127438973097559
$ws.Range("L21").Value = -0.007654384626175384
$ws.Range("M21").Value = 0.4563434365577373
$ws.Range("O21").Value = 0.07258737573549502
$ws.Range("P21").Value = -0.09935659808626392
$ws.Range("C22").Value = 0.008622298328891932
$ws.Range("D22").Value = 0.007140949245637969
$ws.Range("E22").Value = -0.008556948438277938
$ws.Range("G22").Value = -0.01120222345608894
$ws.Range("H22").Value = -0.005456121722244868
$ws.Range("I22").Value = -0.0009385170615406823
$ws.Range("J22").Value = 0.007670824629042182
$ws.Range("K22").Value = -0.003772572438902897
$ws.Range("L22").Value = -0.01204490976179639
$ws.Range("M22").Value = -0.005130972397238895
$ws.Range("O22").Value = -0.01640521784020871
$ws.Range("P22").Value = 0.007994872255794891
$ws.Range("C23").Value = 0.06176445731857828
$ws.Range("D23").Value = 0.003069008474760338
$ws.Range("E23").Value = 0.2762791139631645
$ws.Range("G23").Value = 0.4319510131180405
$ws.Range("H23").Value = 0.007680950131238004
$ws.Range("I23").Value = -0.01020257877610315
$ws.Range("J23").Value = -0.06633927977667689
$ws.Range("K23").Value = -0.0353853202314128
$ws.Range("L23").Value = 0.01430910431636417
$ws.Range("M23").Value = -0.7404759409470375
$ws.Range("O23").Value = -0.1389748927589957
$ws.Range("P23").Value = -0.03310948596437943

# --- Step 4: merge the new Element label cells ---
$ws.Range("A17:A18").Merge() | Out-Null
$ws.Range("A20:A21").Merge() | Out-Null
$ws.Range("A22:A23").Merge() | Out-Null
